$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.444.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.920.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08230"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.010"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.908.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.089"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.241"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06868"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.447.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.669"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.187"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.133.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.573"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.115"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.017"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09654"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.628"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.554"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.378"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02287"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06110"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.180"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.30%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.070"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5963"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.391"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07602"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5593"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.956"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.427"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
